$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 20:40"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1780824
$ws.Range("C4").Value = 12363
$ws.Range("E4").Value = 1174755
$ws.Range("G4").Value = 623
$ws.Range("H4").Value = 103953

# Row 12: India
$ws.Range("B12").Value = 173458
$ws.Range("C12").Value = 8072
$ws.Range("D12").Value = 82468
$ws.Range("E12").Value = 86010
$ws.Range("G12").Value = 269
$ws.Range("H12").Value = 4980

# Row 17: Canada
$ws.Range("D17").Value = 47441
$ws.Range("E17").Value = 34966

# Row 121: Paraguay
$ws.Range("B121").Value = 917
$ws.Range("C121").Value = 17
$ws.Range("D121").Value = 413
$ws.Range("E121").Value = 493
